# Regenerate the lattice-multiplication practice sheet: every cell in the
# 5x3 table keeps its layout (problem line / factor line / rule / two
# partial-product rows) but gets new multiplication facts. The table shape
# (rows/cols/cells) is unchanged, so we just rewrite each cell's Range.Text
# wholesale using a vertical-tab (chr 11) as the in-paragraph line break,
# which is how Word represents <w:br/> inside Range.Text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$cell = $t.Cell(1, 1)
$cell.Range.Text = (@("70 x 23", "  2    3", "  ----", "7|    |", "0|    |") -join $nl)

$cell = $t.Cell(1, 2)
$cell.Range.Text = (@("61 x 66", "  6    6", "  ----", "6|    |", "1|    |") -join $nl)

$cell = $t.Cell(1, 3)
$cell.Range.Text = (@("87 x 26", "  2    6", "  ----", "8|    |", "7|    |") -join $nl)

$cell = $t.Cell(2, 1)
$cell.Range.Text = (@("94 x 79", "  7    9", "  ----", "9|    |", "4|    |") -join $nl)

$cell = $t.Cell(2, 2)
$cell.Range.Text = (@("77 x 72", "  7    2", "  ----", "7|    |", "7|    |") -join $nl)

$cell = $t.Cell(2, 3)
$cell.Range.Text = (@("17 x 95", "  9    5", "  ----", "1|    |", "7|    |") -join $nl)

$cell = $t.Cell(3, 1)
$cell.Range.Text = (@("64 x 50", "  5    0", "  ----", "6|    |", "4|    |") -join $nl)

$cell = $t.Cell(3, 2)
$cell.Range.Text = (@("87 x 84", "  8    4", "  ----", "8|    |", "7|    |") -join $nl)

$cell = $t.Cell(3, 3)
$cell.Range.Text = (@("60 x 83", "  8    3", "  ----", "6|    |", "0|    |") -join $nl)

$cell = $t.Cell(4, 1)
$cell.Range.Text = (@("71 x 87", "  8    7", "  ----", "7|    |", "1|    |") -join $nl)

$cell = $t.Cell(4, 2)
$cell.Range.Text = (@("82 x 57", "  5    7", "  ----", "8|    |", "2|    |") -join $nl)

$cell = $t.Cell(4, 3)
$cell.Range.Text = (@("89 x 63", "  6    3", "  ----", "8|    |", "9|    |") -join $nl)

$cell = $t.Cell(5, 1)
$cell.Range.Text = (@("91 x 49", "  4    9", "  ----", "9|    |", "1|    |") -join $nl)

$cell = $t.Cell(5, 2)
$cell.Range.Text = (@("61 x 41", "  4    1", "  ----", "6|    |", "1|    |") -join $nl)

$cell = $t.Cell(5, 3)
$cell.Range.Text = (@("13 x 94", "  9    4", "  ----", "1|    |", "3|    |") -join $nl)
